# The deck's applied Design ("Integral") and its theme color palette are
# swapped with the plain default "Office Theme" palette that previously
# only backed the (unused-in-COM) secondary theme part. Concretely this
# means driving the presentation's live ColorScheme - reachable from the
# SlideMaster - from the Integral greens/golds to the stock Office blues.
#
# RGB() packs as 0x00BBGGRR (classic OLE COLORREF), so bytes must be
# reversed relative to the RRGGBB hex strings from the target theme.
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

$cs.Item(1).RGB  = (RGBVal 0x00 0x00 0x00)   # dk1
$cs.Item(2).RGB  = (RGBVal 0xFF 0xFF 0xFF)   # lt1
$cs.Item(3).RGB  = (RGBVal 0x44 0x54 0x6A)   # dk2
$cs.Item(4).RGB  = (RGBVal 0xE7 0xE6 0xE6)   # lt2
$cs.Item(5).RGB  = (RGBVal 0x5B 0x9B 0xD5)   # accent1
$cs.Item(6).RGB  = (RGBVal 0xED 0x7D 0x31)   # accent2
$cs.Item(7).RGB  = (RGBVal 0xA5 0xA5 0xA5)   # accent3
$cs.Item(8).RGB  = (RGBVal 0xFF 0xC0 0x00)   # accent4
$cs.Item(9).RGB  = (RGBVal 0x44 0x72 0xC4)   # accent5
$cs.Item(10).RGB = (RGBVal 0x70 0xAD 0x47)   # accent6
$cs.Item(11).RGB = (RGBVal 0x05 0x63 0xC1)   # hlink
$cs.Item(12).RGB = (RGBVal 0x95 0x4F 0x72)   # folHlink
